# Adding day 1 code examples
# Merge the "Class" + " Name" runs (highlighted cyan) into a single
# "ClassName" run on the "Anatomy of a Class" slide.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(13)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

$fullText = $tr.Text
$target = "Class Name"
$startPos = $fullText.IndexOf($target) + 1

$sub = $tr.Characters($startPos, $target.Length)
$sub.Text = "ClassName"
